$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update daily values for column O (day 14) that were previously 0
$ws.Range("O2").Value = 14226.84
$ws.Range("O3").Value = 2253.06
$ws.Range("O4").Value = 2877.5
$ws.Range("O5").Value = 1921.01
$ws.Range("O6").Value = 21278.41

# Update corresponding totals in column AG
$ws.Range("AG2").Value = 120454.79
$ws.Range("AG3").Value = 46967.11
$ws.Range("AG4").Value = 44149.9
$ws.Range("AG5").Value = 35142.86
$ws.Range("AG6").Value = 246714.66
